$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 42890

$ws.Range("B12").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = 1

$ws.Range("C13").Value = "Nasazení na testovací WWW bujabeza, odladění nějakých chyb, náhled pro Jarču"
